# Update Pomc-Mc1r LR-pair sheet with new TPM-derived values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per data row (rows 2-7), for columns G, H, M, N, O, P, Q, R, S, T.
$updates = @{
    2 = @{ G = 1.524170333333333;  H = 4.572511;
           M = 0.7917353333333333; N = 2.375206;
           O = 0.3370298876305008; P = 0.3370298876305008;
           Q = 1.206739506918445;  R = 10.860655562266;
           S = 0.3370298876305008; T = 0.3370298876305008 }
    3 = @{ G = 1.524170333333333;  H = 4.572511;
           M = 0.4973860000000001; N = 1.492158;
           O = 0.2117297796767745; P = 0.2117297796767745;
           Q = 0.7581009854153334; R = 6.822908868738002;
           S = 0.2117297796767745; T = 0.2117297796767745 }
    4 = @{ G = 1.524170333333333;  H = 4.572511;
           M = 0.339559;           N = 1.018677;
           O = 0.1445451867508653; P = 0.1445451867508653;
           Q = 0.5175457542163333; R = 4.657911787947;
           S = 0.1445451867508653; T = 0.1445451867508653 }
    5 = @{ G = 1.524170333333333;  H = 4.572511;
           M = 0.290216;           N = 0.870648;
           O = 0.1235406097853071; P = 0.1235406097853071;
           Q = 0.4423386174586666; R = 3.981047557128;
           S = 0.1235406097853071; T = 0.1235406097853071 }
    6 = @{ G = 1.524170333333333;  H = 4.572511;
           M = 0.116919;           N = 0.350757;
           O = 0.04977066927904846; P = 0.04977066927904845;
           Q = 0.178204471203;     R = 1.603840240827;
           S = 0.04977066927904846; T = 0.04977066927904845 }
    7 = @{ G = 1.524170333333333;  H = 4.572511;
           M = 0.3133393333333334; N = 0.940018;
           O = 0.1333838668775038; P = 0.1333838668775037;
           Q = 0.4775825161331111; R = 4.298242645198001;
           S = 0.1333838668775038; T = 0.1333838668775037 }
}

foreach ($row in $updates.Keys) {
    foreach ($col in $updates[$row].Keys) {
        $ws.Range("$col$row").Value = $updates[$row][$col]
    }
}
